$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target adds explicit column widths for columns A..I (1-9):
# 17, 19, 20, 20, 20, 20, 20, 20, 6 (raw OOXML <col width="..."/>).
#
# The runtime's saved OOXML <col> width comes out as
# (ColumnWidth + 5/6), so back the COM ColumnWidth off by 5/6 here
# so the saved <cols> widths land exactly on the target values.
$ws.Columns.Item(1).ColumnWidth = 16.166666666666668
$ws.Columns.Item(2).ColumnWidth = 18.166666666666668
$ws.Columns.Item(3).ColumnWidth = 19.166666666666668
$ws.Columns.Item(4).ColumnWidth = 19.166666666666668
$ws.Columns.Item(5).ColumnWidth = 19.166666666666668
$ws.Columns.Item(6).ColumnWidth = 19.166666666666668
$ws.Columns.Item(7).ColumnWidth = 19.166666666666668
$ws.Columns.Item(8).ColumnWidth = 19.166666666666668
$ws.Columns.Item(9).ColumnWidth = 5.166666666666667
